$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 125003384
$ws.Range("J64").Value = 4740
$ws.Range("L64").Value = 4740
$ws.Range("N64").Value = -5236

$ws.Range("H67").Value = 125003384
$ws.Range("J67").Value = 4740
$ws.Range("L67").Value = 4740
$ws.Range("N67").Value = -6456

$ws.Range("H74").Value = 3414.923
$ws.Range("I74").Value = 2768.6667
$ws.Range("K74").Value = 2768.6667
$ws.Range("M74").Value = -1832.6667

$ws.Range("H76").Value = 5795.222
$ws.Range("I76").Value = 4668
$ws.Range("J76").Value = 6020.6665
$ws.Range("K76").Value = 4668
$ws.Range("L76").Value = 6020.6665
$ws.Range("M76").Value = -4353
$ws.Range("N76").Value = -6650.6665

$ws.Range("H77").Value = 3414.923
$ws.Range("I77").Value = 2768.6667
$ws.Range("K77").Value = 13843.3335
$ws.Range("M77").Value = -9163.333500000001

$ws.Range("H79").Value = 5795.222
$ws.Range("I79").Value = 4668
$ws.Range("J79").Value = 6020.6665
$ws.Range("K79").Value = 4668
$ws.Range("L79").Value = 6020.6665
$ws.Range("M79").Value = -3576
$ws.Range("N79").Value = -8204.666499999999

$ws.Range("H138").Value = 2869.9297
$ws.Range("I138").Value = 2008.6666
$ws.Range("J138").Value = 3617.8684
$ws.Range("K138").Value = 6025.9998
$ws.Range("L138").Value = 10853.6052
$ws.Range("M138").Value = -885.9997999999996
$ws.Range("N138").Value = -21133.6052

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3686.5715
$ws.Range("I88").Value = 2006
$ws.Range("J88").Value = 3966.6667
$ws.Range("K88").Value = 2006
$ws.Range("L88").Value = 3966.6667
$ws.Range("M88").Value = -1600
$ws.Range("N88").Value = -4778.6667

$ws.Range("H91").Value = 3686.5715
$ws.Range("I91").Value = 2006
$ws.Range("J91").Value = 3966.6667
$ws.Range("K91").Value = 2006
$ws.Range("L91").Value = 3966.6667
$ws.Range("M91").Value = -602
$ws.Range("N91").Value = -6774.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 36724
$ws.Range("J40").Value = 36724
$ws.Range("L40").Value = 36724
$ws.Range("N40").Value = -37254

$ws.Range("H96").Value = 9168
$ws.Range("I96").Value = 6839
$ws.Range("J96").Value = 27800
$ws.Range("K96").Value = 6839
$ws.Range("L96").Value = 27800
$ws.Range("M96").Value = -4093
$ws.Range("N96").Value = -33292

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6027.971
$ws.Range("I132").Value = 5614.1665
$ws.Range("J132").Value = 6930.8184
$ws.Range("K132").Value = 16842.4995
$ws.Range("L132").Value = 20792.4552
$ws.Range("M132").Value = -14312.4995
$ws.Range("N132").Value = -25852.4552

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 58857300
$ws.Range("J141").Value = 71467430
$ws.Range("L141").Value = 71467430
$ws.Range("N141").Value = -71477790

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 953.5333000000001
$ws.Range("I5").Value = 628.55554
$ws.Range("J5").Value = 1441
$ws.Range("K5").Value = 1885.66662
$ws.Range("L5").Value = 4323
$ws.Range("M5").Value = -1773.66662
$ws.Range("N5").Value = -4547

$ws.Range("H87").Value = 5756.5713
$ws.Range("I87").Value = 3806
$ws.Range("J87").Value = 7707.143
$ws.Range("K87").Value = 11418
$ws.Range("L87").Value = 23121.429
$ws.Range("M87").Value = -10170
$ws.Range("N87").Value = -25617.429

$ws.Range("H90").Value = 5756.5713
$ws.Range("I90").Value = 3806
$ws.Range("J90").Value = 7707.143
$ws.Range("K90").Value = 34254
$ws.Range("L90").Value = 69364.287
$ws.Range("M90").Value = -28014
$ws.Range("N90").Value = -81844.287

$ws.Range("H113").Value = 606299.5
$ws.Range("I113").Value = 405.2
$ws.Range("J113").Value = 1688253.6
$ws.Range("K113").Value = 1215.6
$ws.Range("L113").Value = 5064760.800000001
$ws.Range("M113").Value = 954.4000000000001
$ws.Range("N113").Value = -5069100.800000001

$ws.Range("H117").Value = 1131.1111
$ws.Range("I117").Value = 490
$ws.Range("J117").Value = 1314.2858
$ws.Range("K117").Value = 1470
$ws.Range("L117").Value = 3942.8574
$ws.Range("M117").Value = 1972
$ws.Range("N117").Value = -10826.8574

$ws.Range("H135").Value = 953.5333000000001
$ws.Range("I135").Value = 628.55554
$ws.Range("J135").Value = 1441
$ws.Range("K135").Value = 5656.99986
$ws.Range("L135").Value = 12969
$ws.Range("M135").Value = -3121.99986
$ws.Range("N135").Value = -18039

$ws.Range("H141").Value = 2441.611
$ws.Range("I141").Value = 2173.4707
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 6520.4121
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -1340.4121
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19900
$ws.Range("J46").Value = 19900
$ws.Range("L46").Value = 19900
$ws.Range("N46").Value = -20212

$ws.Range("H70").Value = 4015.4285
$ws.Range("I70").Value = 4034.2222
$ws.Range("J70").Value = 3981.6
$ws.Range("K70").Value = 4034.2222
$ws.Range("L70").Value = 3981.6
$ws.Range("M70").Value = -3764.2222
$ws.Range("N70").Value = -4521.6

$ws.Range("H73").Value = 4015.4285
$ws.Range("I73").Value = 4034.2222
$ws.Range("J73").Value = 3981.6
$ws.Range("K73").Value = 4034.2222
$ws.Range("L73").Value = 3981.6
$ws.Range("M73").Value = -3098.2222
$ws.Range("N73").Value = -5853.6

$ws.Range("H80").Value = 3602.5
$ws.Range("I80").Value = 2205
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 2205
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -1207
$ws.Range("N80").Value = -6996

$ws.Range("H83").Value = 3602.5
$ws.Range("I83").Value = 2205
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 11025
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -6033
$ws.Range("N83").Value = -34984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H132").Value = 35182.207
$ws.Range("I132").Value = 6276.48
$ws.Range("J132").Value = 115475.89
$ws.Range("K132").Value = 18829.44
$ws.Range("L132").Value = 346427.67
$ws.Range("M132").Value = -16299.44
$ws.Range("N132").Value = -351487.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 40999.57
$ws.Range("J80").Value = 40999.57
$ws.Range("L80").Value = 40999.57
$ws.Range("N80").Value = -42995.57

$ws.Range("H83").Value = 40999.57
$ws.Range("J83").Value = 40999.57
$ws.Range("L83").Value = 122998.71
$ws.Range("N83").Value = -132982.71

$ws.Range("H132").Value = 5734.971
$ws.Range("I132").Value = 6805.2085
$ws.Range("J132").Value = 3399.9092
$ws.Range("K132").Value = 20415.6255
$ws.Range("L132").Value = 10199.7276
$ws.Range("M132").Value = -17885.6255
$ws.Range("N132").Value = -15259.7276
